# +alter name for mega/diff form
#
# Inserts a new "alterName" column (P) right after the existing "chName"
# column (O), shifting the old "mega" flag column from P to Q. The new
# column starts out as a copy of chName (same Chinese names) and then a
# handful of rows get a genuine alternate name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column O (chName) and insert the copy before column P (mega flag),
# which pushes mega flag to column Q and leaves the new column P with the
# same values/styles as O, exactly like the source workbook's edit.
$ws.Columns.Item(15).Copy()
$ws.Columns.Item(16).Insert()

# Header for the new column.
$ws.Range("P1").Value = "alterName"

# Rows whose alternate name differs from chName (order matches the
# original author's edit order so the shared-strings table lines up too).
$ws.Range("P11").Value = "超级袋龙"
$ws.Range("P12").Value = "超级大甲"
$ws.Range("P22").Value = "超级班吉拉"
$ws.Range("P36").Value = "超级灾兽"
$ws.Range("P38").Value = "超级血翼飞龙"
$ws.Range("P39").Value = "超级合金十字"
$ws.Range("P44").Value = "超级裂空座"
$ws.Range("P7").Value = "超级比雕"

# Best-effort view-state cosmetics matching the author's edit.
$win = $wb.Windows.Item(1)
$win.Left = 10440

$ws.Range("P46").Select()
